$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New phase-distribution values for A2:B107 (rows 2..107), replacing the
# previous recursive/matching-layer loss curves with the reversed RT data.
$data = @(
    @(-0.5177348774463787,-482.8672876447307),
    @(-0.5078098721269414,-481.6889470267456),
    @(-0.4978848667350946,-480.5106064001636),
    @(-0.4879598613794522,-479.3322657778801),
    @(-0.4780348560592672,-478.1539251598061),
    @(-0.4681098506681684,-476.9755845333129),
    @(-0.4581848453125258,-475.7972439110295),
    @(-0.4482598399568842,-474.6189032887459),
    @(-0.4383348346356992,-473.4405626705533),
    @(-0.4284098292203118,-472.2622220411765),
    @(-0.4184848239239157,-471.0838814259268),
    @(-0.4085598185343151,-469.9055407996117),
    @(-0.3986348132121328,-468.7272001813007),
    @(-0.3887098078230316,-467.5488595550447),
    @(-0.378784802443225,-466.3705189298922),
    @(-0.3688597971117475,-465.1921783104776),
    @(-0.358934791764849,-464.0138376892321),
    @(-0.3490097864004638,-462.8354970659106),
    @(-0.3390847810448211,-461.6571564436269),
    @(-0.3291597757208907,-460.4788158251083),
    @(-0.3192347703335372,-459.3004751990599),
    @(-0.3093097650091061,-458.1221345804819),
    @(-0.2993847596532149,-456.9437939581687),
    @(-0.2894597542444678,-455.7654533295803),
    @(-0.2795347489414319,-454.5871127135425),
    @(-0.2696097435855395,-453.4087720912291),
    @(-0.2596847381996842,-452.2304314653586),
    @(-0.2497597328440431,-451.0520908430752),
    @(-0.2398347275178642,-449.8737502242897),
    @(-0.229909722132759,-448.6954095985081),
    @(-0.2199847167771163,-447.5170689762245),
    @(-0.2100597114214736,-446.3387283539408),
    @(-0.2001347060942971,-445.1603877350369),
    @(-0.1902097007384047,-443.9820471127236),
    @(-0.180284695382514,-442.8037064904105),
    @(-0.1703596900266216,-441.6253658680974),
    @(-0.1604346846432633,-440.4470252425232),
    @(-0.1505096793148384,-439.2686846234709),
    @(-0.140584673958946,-438.0903440011577),
    @(-0.1306596686030553,-436.9120033788447),
    @(-0.1207346632206957,-435.7336627533891),
    @(-0.1108096578650527,-434.5553221311055),
    @(-0.1008846525094117,-433.3769815088221),
    @(-0.09095964717948724,-432.1986408895918),
    @(-0.0810346418061725,-431.0203002652101),
    @(-0.07110963642528623,-429.8419596399294),
    @(-0.06118463108684238,-428.6636190196878),
    @(-0.05125962575592068,-427.4852784003392),
    @(-0.04133462040002839,-426.3069377780259),
    @(-0.03140961504413733,-425.1285971557128),
    @(-0.02148460968824496,-423.9502565333996),
    @(-0.0115596043086321,-422.77191590827),
    @(-0.001634598976461639,-421.5935752887731),
    @(0.008290406387131301,-420.4152346655457),
    @(0.01821541173532145,-419.2368940441469),
    @(0.02814041709121383,-418.0585534218336),
    @(0.03806542246957818,-416.8802127968524),
    @(0.0479904278029971,-415.7018721772073),
    @(0.05791543315888949,-414.5235315548941),
    @(0.06784043850072774,-413.3451909342493),
    @(0.07776544389214612,-412.1668503077183),
    @(0.0876904492401687,-410.9885096863394),
    @(0.0976154546034298,-409.8101690631513),
    @(0.1075404599449252,-408.6318284425473),
    @(0.1174654652942393,-407.453487821015),
    @(0.1273904706640041,-406.2751471970547),
    @(0.137315476032239,-405.0968065732761),
    @(0.1472404813816405,-403.9184659517335),
    @(0.1571654867312679,-402.740125330164),
    @(0.1670904920929241,-401.5617847071665),
    @(0.1770154974427763,-400.3834440855703),
    @(0.1869405027985312,-399.2051034632734),
    @(0.196865508154285,-398.0267628409767),
    @(0.20679051351004,-396.8484222186797),
    @(0.2167155188711345,-395.6700815957488),
    @(0.2266405242215485,-394.491740974086),
    @(0.2365655295773036,-393.313400351789),
    @(0.2464905349330572,-392.1350597294922),
    @(0.2564155402888121,-390.9567191071953),
    @(0.2663405456278367,-389.7783784868847),
    @(0.2762655510003207,-388.6000378626015),
    @(0.2861905563398462,-387.4216972422315),
    @(0.2961155617002915,-386.2433566193777),
    @(0.3060405670675844,-385.0650159957109),
    @(0.3159655724233379,-383.8866753734142),
    @(0.3258905777638621,-382.7083347529255),
    @(0.3358155831388407,-381.5299941283462),
    @(0.3457405884906016,-380.3516535065235),
    @(0.3556655938356427,-379.1733128854985),
    @(0.3655905991878796,-377.9949722636192),
    @(0.3755156045578651,-376.8166316396328),
    @(0.3854406098998875,-375.6382910189662),
    @(0.3953656152558907,-374.4599503966398),
    @(0.4052906206118955,-373.2816097743132),
    @(0.415215625970994,-372.1032691516193),
    @(0.425140631326885,-370.9249285293062),
    @(0.4350656366770391,-369.7465879076742),
    @(0.4449906420331547,-368.5682472853345),
    @(0.4549156473945606,-367.3899066623666),
    @(0.4648406527555168,-366.2115660394521),
    @(0.4747656581111831,-365.0332254171658),
    @(0.4846906634668506,-363.8548847948791),
    @(0.4946156688181271,-362.6765441731139),
    @(0.504540674188921,-361.4982035490315),
    @(0.514465679544339,-360.3198629267746),
    @(0.5243906848997554,-359.1415223045179)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Rows 93:107 are brand new; give column A the same bold/bordered/centered
# style used by the rest of column A (copy format from A2).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A93:A107").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
